# "Generate Report for Handback" -- refresh the localization-status report
# after a handback: status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the per-language handback timestamps
# advance, and the stale "handback file is not latest" error clears now
# that everything is in sync.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-18 14:53:49"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-18 14:53:56"
$wsDeDe.Range("P2").Value = ""

# --- Column widths: the Status column widened to fit the longer new
# status text, and the now-empty Error Detail column shrank back down.
# (ColumnWidth only lands on 1/6-character increments, so we pick the
# closest achievable width to the recorded target.)
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333332

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333332
